$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.250.83"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "3.390.79"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "594.69"
$ws.Range("E5").Value = "  +7.27%  "
$ws.Range("D6").Value = "188.01"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.603"
$ws.Range("E7").Value = "  +4.06%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("D10").Value = "0.592"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "47.83"
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").Value = "  +5.64%  "
$ws.Range("D13").Value = "3.925.04"
$ws.Range("D14").Value = "641.06"
$ws.Range("E14").Value = "  +10.51%  "
$ws.Range("D15").Value = "8.66"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "68.055.43"
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").Value = "3.382.96"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "18.14"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "11.14"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "0.915"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "17.95"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").Value = "100.18"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("E26").Value = "  +6.33%  "
$ws.Range("E27").Value = "  +4.43%  "
$ws.Range("D28").Value = "32.98"
$ws.Range("E28").Value = "  +7.85%  "
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  +4.63%  "
$ws.Range("D31").Value = "615.79"
$ws.Range("E31").Value = "  +6.80%  "
$ws.Range("D32").Value = "3.88"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "4.039.01"
$ws.Range("E33").Value = "  +8.62%  "
$ws.Range("D34").Value = "11.16"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "56.42"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").Value = "2.80"
$ws.Range("E38").Value = "  +6.91%  "
$ws.Range("E39").Value = "  +3.62%  "
$ws.Range("D40").Value = "33.88"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "3.29"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").Value = "0.0₃0710"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "3.42"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.345"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").Value = "0.0425"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("D46").Value = "0.131"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "2.62"
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("E48").Value = "  +12.73%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "128.88"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  +6.82%  "
